$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Usuarios")
$ws2 = $wb.Worksheets.Item("Recursos")

# ---------------------------------------------------------------------------
# Hoja "Usuarios": deja solo la cabecera y una unica fila de datos
# (antes habia 13 usuarios individuales, ahora solo "Perez Gonzalez, Docente")
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value2 = "Pérez González, Docente"
$null = $ws1.Range("A3:A14").EntireRow.Delete()

$ws1.Columns("A").ColumnWidth = 28.75

$null = $ws1.Range("A3:A11").Select()

# ---------------------------------------------------------------------------
# Hoja "Recursos": se inserta una nueva columna "IDRecurso" entre
# "Contexto del evento" y "Alias", y se actualizan los datos de ejemplo.
# ---------------------------------------------------------------------------
$null = $ws2.Columns("B").Insert()

$ws2.Range("B1").Value2 = "IDRecurso"

$ws2.Range("A2").Value2 = "Foro: Noticias de clase"
$ws2.Range("B2").Value2 = 5000
$ws2.Range("C2").Value2 = "Foro: Noticias de clase"

$ws2.Columns("A").ColumnWidth = 29.75
$ws2.Columns("B").ColumnWidth = 25.1
$ws2.Columns("C").ColumnWidth = 29.75

$null = $ws2.Activate()
$null = $ws2.Range("C2").Select()
